$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2020" data column (K), mirroring the existing "2019" column (J)
# cell-by-cell so its formatting (border/number style) matches column J.
# xlPasteFormats = -4122
$xlPasteFormats = -4122

# K3: bottom-border divider cell above the year header row (no value).
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial($xlPasteFormats)

# K4: year header "2020".
$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial($xlPasteFormats)
$ws.Range("K4").Value = 2020

# K6: "Mammals" 2020 figure.
$ws.Range("J6").Copy()
$ws.Range("K6").PasteSpecial($xlPasteFormats)
$ws.Range("K6").Value = 5.9

# K7: "Birds" 2020 figure.
$ws.Range("J7").Copy()
$ws.Range("K7").PasteSpecial($xlPasteFormats)
$ws.Range("K7").Value = 1.5

# K8: "Amphibians and Reptiles" 2020 figure - no data available ("-").
$ws.Range("J8").Copy()
$ws.Range("K8").PasteSpecial($xlPasteFormats)
$ws.Range("K8").Value = "-"

# Matches the authored selection state recorded in the saved workbook.
$ws.Range("L16").Select()
